$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.302.08'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').Value = '2.493.16'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('D5').Value = '321.06'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').Value = '108.47'
$ws.Range('E6').Value = '  +3.82%  '
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '0.535'
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('D10').Value = '39.08'
$ws.Range('E10').Value = '  +5.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0810'
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').Value = '18.39'
$ws.Range('E13').Value = '  +0.74%  '
$ws.Range('D14').Value = '7.14'
$ws.Range('E14').Value = '  -0.45%  '
$ws.Range('D15').Value = '2.880.72'
$ws.Range('E15').Value = '  +0.30%  '
$ws.Range('D16').Value = '2.488.54'
$ws.Range('E16').Value = '  -1.76%  '
$ws.Range('D17').Value = '0.847'
$ws.Range('E17').Value = '  +0.81%  '
$ws.Range('D18').Value = '47.204.06'
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('D19').Value = '13.01'
$ws.Range('E19').Value = '  +3.66%  '
$ws.Range('D20').Value = '6.61'
$ws.Range('E20').Value = '  +0.60%  '
$ws.Range('D21').Value = '0.0₃0934'
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('D22').Value = '2.67'
$ws.Range('E22').Value = '  +13.55%  '
$ws.Range('D23').Value = '70.35'
$ws.Range('E23').Value = '  -0.24%  '
$ws.Range('D24').Value = '245.45'
$ws.Range('E24').Value = '  -2.00%  '
$ws.Range('D25').Value = '2.56'
$ws.Range('E25').Value = '  +0.51%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').Value = '25.77'
$ws.Range('E27').Value = '  -1.15%  '
$ws.Range('E28').Value = '  +3.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.00'
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = '34.73'
$ws.Range('E30').Value = '  -0.93%  '
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').Value = '0.137'
$ws.Range('E31').Value = '  +3.37%  '
$ws.Range('D32').Value = '49.69'
$ws.Range('E32').Value = '  +0.67%  '
$ws.Range('D33').Value = '20.64'
$ws.Range('E33').Value = '  +5.72%  '
$ws.Range('E34').Value = '  +1.05%  '
$ws.Range('D35').Value = '0.0783'
$ws.Range('E35').Value = '  +1.28%  '
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('D37').Value = '4.76'
$ws.Range('E37').Value = '  +4.37%  '
$ws.Range('D38').Value = '1.97'
$ws.Range('E38').Value = '  +2.49%  '
$ws.Range('E39').Value = '  -0.44%  '
$ws.Range('D40').Value = '23.03'
$ws.Range('E40').Value = '  +7.61%  '
$ws.Range('E41').Value = '  -0.24%  '
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('D43').Value = '117.28'
$ws.Range('E43').Value = '  -4.02%  '
$ws.Range('E44').Value = '  +0.38%  '
$ws.Range('D45').Value = '1.994.00'
$ws.Range('E45').Value = '  +2.05%  '
$ws.Range('D46').Value = '3.04'
$ws.Range('E46').Value = '  +2.48%  '
$ws.Range('D47').Value = '2.01'
$ws.Range('E47').Value = '  -5.16%  '
$ws.Range('D48').Value = '9.13'
$ws.Range('E48').Value = '  +0.12%  '
$ws.Range('D49').Value = '1.78'
$ws.Range('E49').Value = '  -0.48%  '
$ws.Range('E50').Value = '  -5.39%  '
$ws.Range('D51').Value = '56.64'
$ws.Range('E51').Value = '  +4.26%  '
